# Append " (Changed main)" after the existing sentence in the first
# paragraph, as three distinct runs:
#   <w:r><w:t xml:space="preserve"> (</w:t></w:r>
#   <w:r><w:t>Changed main</w:t></w:r>
#   <w:r><w:t>)</w:t></w:r>
#
# A plain Range.InsertAfter() would grow/merge into the existing run
# (same, empty run formatting), collapsing everything back down to a
# single <w:r>. To preserve the three separate runs exactly as authored,
# splice the new run XML straight into the paragraph's OOXML and push it
# back in with Range.InsertXML (which replaces - rather than appends to -
# the target range's contents).

$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range

$full = $r1.WordOpenXML

# Namespace/attribute declarations carried on <w:document ...> so the
# qualified attribute names (w14:paraId etc.) on <w:p> keep resolving.
if ($full -match '(?s)<w:document([^>]*)>') {
    $docAttrs = $matches[1]
} else {
    $docAttrs = ''
}

# The first paragraph's current OOXML, exactly as it stands today.
if ($full -match '(?s)<w:body>(<w:p\b.*?</w:p>)') {
    $paraXml = $matches[1]
} else {
    throw "edit.ps1: could not locate first paragraph in WordOpenXML"
}

$newRuns = '<w:r><w:t xml:space="preserve"> (</w:t></w:r>' +
           '<w:r><w:t>Changed main</w:t></w:r>' +
           '<w:r><w:t>)</w:t></w:r>'

$updatedParaXml = $paraXml -replace '</w:p>$', ($newRuns + '</w:p>')

$wrapped = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document' + $docAttrs + '>' +
           '<w:body>' + $updatedParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r1.InsertXML($wrapped) | Out-Null
